# Apply updated odds values to Sheet1, matching the upstream FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("L4").Value = 1.5
$ws.Range("M4").Value = 2.63
$ws.Range("N4").Value = 2.5
$ws.Range("O4").Value = 1.53

# Row 6
$ws.Range("G6").Value = 3.4
$ws.Range("I6").Value = 2.3
$ws.Range("T6").Value = 7
$ws.Range("V6").Value = 13
$ws.Range("X6").Value = 34
$ws.Range("Z6").Value = 6
$ws.Range("AF6").Value = 9.5
$ws.Range("AG6").Value = 11
$ws.Range("AH6").Value = 23

# Row 10
$ws.Range("G10").Value = 1.8
$ws.Range("I10").Value = 4.2
$ws.Range("R10").Value = 1.67
$ws.Range("S10").Value = 2.1
$ws.Range("T10").Value = 8.5
$ws.Range("U10").Value = 9.5
$ws.Range("Y10").Value = 23
$ws.Range("AI10").Value = 34

# Row 13
$ws.Range("G13").Value = 2.1
$ws.Range("I13").Value = 3.55
$ws.Range("K13").Value = 6.2
$ws.Range("U13").Value = 9.25
$ws.Range("Z13").Value = 6.2
$ws.Range("AA13").Value = 6.6
$ws.Range("AE13").Value = 8
$ws.Range("AF13").Value = 18
$ws.Range("AI13").Value = 45

# Row 17
$ws.Range("G17").Value = 2.55
$ws.Range("H17").Value = 2.88
$ws.Range("I17").Value = 3.1
$ws.Range("L17").Value = 1.5
$ws.Range("M17").Value = 2.5
$ws.Range("T17").Value = 6.5
$ws.Range("U17").Value = 11
$ws.Range("V17").Value = 11
$ws.Range("W17").Value = 26
$ws.Range("X17").Value = 23
$ws.Range("AE17").Value = 7.5
$ws.Range("AF17").Value = 13
$ws.Range("AG17").Value = 12
$ws.Range("AH17").Value = 34
$ws.Range("AI17").Value = 29

# Row 19
$ws.Range("G19").Value = 2.8
$ws.Range("H19").Value = 2.85
$ws.Range("I19").Value = 2.62
$ws.Range("N19").Value = 2.4
$ws.Range("P19").Value = 1.57
$ws.Range("Q19").Value = 2.12
$ws.Range("R19").Value = 1.98
$ws.Range("T19").Value = 6.9
$ws.Range("U19").Value = 13
$ws.Range("V19").Value = 10.75
$ws.Range("W19").Value = 35
$ws.Range("X19").Value = 28
$ws.Range("Y19").Value = 45
$ws.Range("Z19").Value = 6.3
$ws.Range("AA19").Value = 5.7
$ws.Range("AB19").Value = 17
$ws.Range("AE19").Value = 6.3
$ws.Range("AF19").Value = 11.5
$ws.Range("AG19").Value = 10.5
$ws.Range("AH19").Value = 30
$ws.Range("AI19").Value = 28
$ws.Range("AJ19").Value = 45

# Row 20
$ws.Range("G20").Value = 1.53
$ws.Range("H20").Value = 3.8
$ws.Range("I20").Value = 5.7
$ws.Range("L20").Value = 1.33
$ws.Range("M20").Value = 2.77
$ws.Range("N20").Value = 1.98
$ws.Range("P20").Value = 1.42
$ws.Range("Q20").Value = 2.45
$ws.Range("R20").Value = 2.07
$ws.Range("S20").Value = 1.6
$ws.Range("T20").Value = 5.6
$ws.Range("U20").Value = 6.3
$ws.Range("W20").Value = 10.25
$ws.Range("X20").Value = 14
$ws.Range("Z20").Value = 9
$ws.Range("AA20").Value = 7.6
$ws.Range("AB20").Value = 22
$ws.Range("AE20").Value = 12.5
$ws.Range("AF20").Value = 32
$ws.Range("AG20").Value = 19
$ws.Range("AH20").Value = 120
$ws.Range("AI20").Value = 70
$ws.Range("AJ20").Value = 80

# Row 21
$ws.Range("H21").Value = 3.25
$ws.Range("I21").Value = 2.3
$ws.Range("J21").Value = 1.07
$ws.Range("K21").Value = 9
$ws.Range("L21").Value = 1.33
$ws.Range("M21").Value = 3.25
$ws.Range("N21").Value = 2.05
$ws.Range("O21").Value = 1.75
$ws.Range("P21").Value = 1.44
$ws.Range("Q21").Value = 2.63
$ws.Range("Z21").Value = 9
$ws.Range("AA21").Value = 6
$ws.Range("AE21").Value = 7.5

# Row 22
$ws.Range("G22").Value = 1.8
$ws.Range("I22").Value = 4.5
$ws.Range("J22").Value = 1.08
$ws.Range("K22").Value = 8
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 1.73
$ws.Range("U22").Value = 7.5
$ws.Range("Z22").Value = 8
$ws.Range("AE22").Value = 11
$ws.Range("AJ22").Value = 51

# Row 25
$ws.Range("P25").Value = 1.44
$ws.Range("Q25").Value = 2.63

# Row 60
$ws.Range("L60").Value = 1.17
$ws.Range("M60").Value = 5
$ws.Range("N60").Value = 1.6
$ws.Range("O60").Value = 2.3
$ws.Range("R60").Value = 1.47

# Row 61
$ws.Range("G61").Value = 1.65
$ws.Range("H61").Value = 3.9
$ws.Range("I61").Value = 4.45
$ws.Range("L61").Value = 1.17
$ws.Range("U61").Value = 9.25
$ws.Range("Y61").Value = 19.5
$ws.Range("Z61").Value = 14
$ws.Range("AA61").Value = 7.9
$ws.Range("AE61").Value = 15.5
$ws.Range("AF61").Value = 29
$ws.Range("AG61").Value = 14.5
$ws.Range("AH61").Value = 75

# Row 67
$ws.Range("O67").Value = 1.85

# Row 86
$ws.Range("G86").Value = 2.02
$ws.Range("I86").Value = 3.5
$ws.Range("Q86").Value = 2.52
$ws.Range("T86").Value = 7.7
$ws.Range("U86").Value = 10
$ws.Range("W86").Value = 19
$ws.Range("Z86").Value = 9.75
$ws.Range("AA86").Value = 6.3
$ws.Range("AE86").Value = 10.5
$ws.Range("AF86").Value = 19.5
$ws.Range("AG86").Value = 11.75
$ws.Range("AH86").Value = 50
$ws.Range("AI86").Value = 32
